$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.8423696626666437, 0.5676861504354533, 1.3446534841860367, 2.885134191327045)"
$ws.Range("C2").Value = "NIG(1.2683718253433873, 1.0636170667672702, 3.3685861174046763, 4.521407957255073)"
$ws.Range("D2").Value = "NIG(0.8812254748756672, 0.6739984555960135, 1.15119180269309, 3.4136644500348234)"
$ws.Range("E2").Value = "JSU(-1.0201909157321887, 1.3136299823044724, 4.103530009428092, 6.590021246927642)"
